# Delete row 8 (the "apgd2" row with Adv_param 0.2) from the "breastcancer" sheet.
# Excel will automatically shift all subsequent rows up by one, adjust the
# merged cell ranges that spanned the deleted row, and update the used
# range/dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("breastcancer")

$ws.Rows.Item(8).Delete()

# The "apgd2" label lived only in the merged range's anchor cell (A8), which
# was removed along with the deleted row. Restore it on the new anchor cell
# (the row that shifted up into position 8).
$ws.Range("A8").Value = "apgd2"
